$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019-S1")

# Row 3 (Spring Cloud): ring moves from "trial" to "adopt"
$ws.Range("B3").Value = "adopt"

# Row 6: MongoDB -> Metrics by Dropwizard, add description, row height 29
$ws.Range("A6").Value = "Metrics by Dropwizard"
$ws.Range("E6").Value = "Metrics provides a powerful toolkit of ways to measure the behavior of critical components in your production environment"
$ws.Rows.Item(6).RowHeight = 29

# Row 7: Google Cloud Platform -> Apache Kafka, ring trial->adopt, quadrant platforms->tools, add description, row height 232
$ws.Range("A7").Value = "Apache Kafka"
$ws.Range("B7").Value = "adopt"
$ws.Range("C7").Value = "tools"
$kafkaDescription = "Kafka is used for building real-time data pipelines and streaming apps. It is horizontally scalable, fault-tolerant, fast, and runs in production in thousands of companies.`nA streaming platform has three key capabilities:`nPublish and subscribe to streams of records, similar to a message queue or enterprise messaging system.`nStore streams of records in a fault-tolerant durable way.`nProcess streams of records as they occur.`nKafka is generally used for two broad classes of applications:`nBuilding real-time streaming data pipelines that reliably get data between systems or applications`nBuilding real-time streaming applications that transform or react to the streams of data"
$ws.Range("E7").Value = $kafkaDescription
$ws.Range("E7").WrapText = $true
$ws.Range("E7").VerticalAlignment = -4160
$ws.Rows.Item(7).RowHeight = 232

# Row 8: Microservices -> Elasticsearch, quadrant techniques->tools, add description
$ws.Range("A8").Value = "Elasticsearch"
$ws.Range("C8").Value = "tools"
$ws.Range("E8").Value = "Elasticsearch is a distributed, RESTful search and analytics engine."

# Update sheet view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A9").Select()
